# Apply updated cryptocurrency price/volume data to Sheet1 (cells D2:E51, with a
# 3-row re-rank of B/C/D/E in rows 32-34), matching the upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Prefix with a literal apostrophe so Excel stores the value as text even when
    # it looks numeric (e.g. "1.010" or "0.3952"), then drop the style back to
    # "Normal" so the quote-prefix formatting flag does not linger on the cell.
    $range.Value = "`'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '24.858.22'
Set-TextCell $ws.Range("E2") '  +0.68%  '
# Row 3
Set-TextCell $ws.Range("D3") '1.708.00'
Set-TextCell $ws.Range("E3") '  +0.79%  '
# Row 4
Set-TextCell $ws.Range("D4") '1.010'
Set-TextCell $ws.Range("E4") '  +0.95%  '
# Row 5
Set-TextCell $ws.Range("D5") '317.45'
Set-TextCell $ws.Range("E5") '  +0.19%  '
# Row 6
Set-TextCell $ws.Range("D6") '1.010'
Set-TextCell $ws.Range("E6") '  +1.01%  '
# Row 7
Set-TextCell $ws.Range("D7") '0.3952'
Set-TextCell $ws.Range("E7") '  -0.10%  '
# Row 8
Set-TextCell $ws.Range("D8") '0.4098'
Set-TextCell $ws.Range("E8") '  +0.71%  '
# Row 9
Set-TextCell $ws.Range("D9") '1.523'
Set-TextCell $ws.Range("E9") '  +2.13%  '
# Row 10
Set-TextCell $ws.Range("D10") '1.010'
Set-TextCell $ws.Range("E10") '  +0.91%  '
# Row 11
Set-TextCell $ws.Range("D11") '52.42'
Set-TextCell $ws.Range("E11") '  +1.24%  '
# Row 12
Set-TextCell $ws.Range("D12") '0.08832'
Set-TextCell $ws.Range("E12") '  -0.76%  '
# Row 13
Set-TextCell $ws.Range("D13") '7.684'
Set-TextCell $ws.Range("E13") '  +6.85%  '
# Row 14
Set-TextCell $ws.Range("D14") '24.74'
Set-TextCell $ws.Range("E14") '  +5.02%  '
# Row 15
Set-TextCell $ws.Range("D15") '0.00001371'
Set-TextCell $ws.Range("E15") '  +3.47%  '
# Row 16
Set-TextCell $ws.Range("D16") '8.075'
Set-TextCell $ws.Range("E16") '  -1.30%  '
# Row 17
Set-TextCell $ws.Range("D17") '1.712.23'
Set-TextCell $ws.Range("E17") '  +1.07%  '
# Row 18
Set-TextCell $ws.Range("D18") '99.96'
Set-TextCell $ws.Range("E18") '  -0.26%  '
# Row 19
Set-TextCell $ws.Range("D19") '0.07119'
Set-TextCell $ws.Range("E19") '  +1.62%  '
# Row 20
Set-TextCell $ws.Range("D20") '20.02'
Set-TextCell $ws.Range("E20") '  +1.42%  '
# Row 21
Set-TextCell $ws.Range("D21") '7.438'
Set-TextCell $ws.Range("E21") '  +5.86%  '
# Row 22
Set-TextCell $ws.Range("E22") '  +1.30%  '
# Row 23
Set-TextCell $ws.Range("E23") '  +0.74%  '
# Row 24
Set-TextCell $ws.Range("D24") '24.856.56'
Set-TextCell $ws.Range("E24") '  +0.65%  '
# Row 25
Set-TextCell $ws.Range("D25") '3.065'
Set-TextCell $ws.Range("E25") '  -3.98%  '
# Row 26
Set-TextCell $ws.Range("E26") '  +0.55%  '
# Row 27
Set-TextCell $ws.Range("D27") '22.93'
Set-TextCell $ws.Range("E27") '  +1.08%  '
# Row 28
Set-TextCell $ws.Range("D28") '164.71'
Set-TextCell $ws.Range("E28") '  +0.83%  '
# Row 29
Set-TextCell $ws.Range("D29") '8.782'
Set-TextCell $ws.Range("E29") '  +15.41%  '
# Row 30
Set-TextCell $ws.Range("D30") '139.05'
Set-TextCell $ws.Range("E30") '  +1.51%  '
# Row 31
Set-TextCell $ws.Range("D31") '5.214'
Set-TextCell $ws.Range("E31") '  +0.97%  '
# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range("D32") '7.725'
Set-TextCell $ws.Range("E32") '  +7.58%  '
# Row 33
$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws.Range("D33") '1.899.45'
Set-TextCell $ws.Range("E33") '  +0.36%  '
# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws.Range("D34") '0.08998'
Set-TextCell $ws.Range("E34") '  +4.51%  '
# Row 35
Set-TextCell $ws.Range("D35") '1.058'
Set-TextCell $ws.Range("E35") '  -1.30%  '
# Row 36
Set-TextCell $ws.Range("D36") '1.986'
Set-TextCell $ws.Range("E36") '  +3.20%  '
# Row 37
Set-TextCell $ws.Range("D37") '0.02945'
Set-TextCell $ws.Range("E37") '  +8.05%  '
# Row 38
Set-TextCell $ws.Range("D38") '0.2757'
Set-TextCell $ws.Range("E38") '  +0.30%  '
# Row 39
Set-TextCell $ws.Range("D39") '10.97'
Set-TextCell $ws.Range("E39") '  -4.35%  '
# Row 40
Set-TextCell $ws.Range("D40") '14.47'
Set-TextCell $ws.Range("E40") '  -0.18%  '
# Row 41
Set-TextCell $ws.Range("D41") '0.09199'
Set-TextCell $ws.Range("E41") '  +0.10%  '
# Row 42
Set-TextCell $ws.Range("D42") '0.7922'
Set-TextCell $ws.Range("E42") '  +3.16%  '
# Row 43
Set-TextCell $ws.Range("D43") '1.475'
Set-TextCell $ws.Range("E43") '  -0.14%  '
# Row 44
Set-TextCell $ws.Range("D44") '16.62'
Set-TextCell $ws.Range("E44") '  +2.53%  '
# Row 45
Set-TextCell $ws.Range("D45") '0.7291'
Set-TextCell $ws.Range("E45") '  +1.33%  '
# Row 46
Set-TextCell $ws.Range("D46") '2.627'
Set-TextCell $ws.Range("E46") '  +0.37%  '
# Row 47
Set-TextCell $ws.Range("D47") '4.268'
Set-TextCell $ws.Range("E47") '  +0.94%  '
# Row 48
Set-TextCell $ws.Range("D48") '1.010'
Set-TextCell $ws.Range("E48") '  +0.94%  '
# Row 49
Set-TextCell $ws.Range("E49") '  +0.45%  '
# Row 50
Set-TextCell $ws.Range("D50") '140.18'
Set-TextCell $ws.Range("E50") '  -0.34%  '
# Row 51
Set-TextCell $ws.Range("D51") '91.93'
Set-TextCell $ws.Range("E51") '  +1.90%  '
